# Delete row 6 (the "Test1" soldier entry) from the active worksheet.
# This mirrors a user selecting the whole row 6 in Excel and choosing
# Delete, which shifts rows 7 and 8 up to become rows 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row first so the saved sheet view reflects the
# same selection a human would have left behind after deleting it
# (activeCell A6, selection spanning the whole row A6:XFD6).
$ws.Rows(6).Select()

# Remove the row; Excel shifts everything below it upward.
$ws.Rows(6).Delete()

# Re-enter the formula across the two rows that used to hold the
# shared formula group (old Q7:Q8) so it stays a single shared
# formula now anchored at Q6:Q7, matching native Excel behavior.
$ws.Range("Q6:Q7").Formula = "=(E6*2)+(F6*4)+G6+H6+I6+(J6*2)"
